# CIERRE 11 JUN 22
# Move the payroll receipt forward one week: "SEMANA 22" (30 May - 04 Jun 2022)
# becomes "SEMANA 23" (06 - 12 Jun 2022), update the advance/discount figures
# for the new week, and re-point the saved selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Week-label text (shared string). B9 is the master cell; H9, B27, H27 and
#    B43 all hold formulas (=B9, =B27, =H27) so they pick the new text up
#    automatically on recalculation.
$ws.Range("B9").Value = "SEMANA   23  DEL    06      Al   12   DE   JUNIO          2022"

# 2) Figures that changed for the new week. These are literal numbers (not
#    formulas); the SUM() totals that depend on them (K24, E41) recompute
#    automatically.
$ws.Range("K21").Value = 1400
$ws.Range("E40").Value = 1250

# 3) Restore the saved selection for this edit (H39:I40, anchored at H39) and
#    scroll position (top-left corner around row 22).
[void]$ws.Range("H39:I40").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
